$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 211
$ws1.Range("F3").Value = 118
$ws1.Range("F4").Value = 407
$ws1.Range("F5").Value = 996
$ws1.Range("F6").Value = 5547
$ws1.Range("F7").Value = 495
$ws1.Range("F8").Value = 698
$ws1.Range("F9").Value = 957
$ws1.Range("F10").Value = 825
$ws1.Range("F11").Value = 80
$ws1.Range("F15").Value = 22
$ws1.Range("F17").Value = 1859
$ws1.Range("F18").Value = 1476
$ws1.Range("F19").Value = 927
$ws1.Range("F21").Value = 196
$ws1.Range("F22").Value = 338
$ws1.Range("F23").Value = 557
$ws1.Range("F24").Value = 157
$ws1.Range("F28").Value = 2957
$ws1.Range("F30").Value = 103
$ws1.Range("F31").Value = 64
$ws1.Range("F32").Value = 124
$ws1.Range("F33").Value = 36
$ws1.Range("F34").Value = 394
$ws1.Range("F37").Value = 14
$ws1.Range("F39").Value = 294
$ws1.Range("F40").Value = 728
$ws1.Range("F42").Value = 52
$ws1.Range("F43").Value = 57
$ws1.Range("F44").Value = 70

# Sheet 2: 演出 (Performance)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F4").Value = 195
$ws2.Range("F6").Value = 134

# Sheet 4: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 211
$ws4.Range("F4").Value = 118
$ws4.Range("F5").Value = 996
$ws4.Range("F7").Value = 5547
$ws4.Range("F8").Value = 495
$ws4.Range("F9").Value = 698
$ws4.Range("F11").Value = 195
$ws4.Range("F12").Value = 957
$ws4.Range("F13").Value = 825
$ws4.Range("F15").Value = 134
$ws4.Range("F16").Value = 80
$ws4.Range("F20").Value = 22
$ws4.Range("F23").Value = 1859
$ws4.Range("F24").Value = 1476
$ws4.Range("F25").Value = 927
$ws4.Range("F26").Value = 196
$ws4.Range("F27").Value = 338
$ws4.Range("F29").Value = 557
$ws4.Range("F30").Value = 157
$ws4.Range("F32").Value = 2958
$ws4.Range("F34").Value = 103
$ws4.Range("F35").Value = 64
$ws4.Range("F36").Value = 124
$ws4.Range("F37").Value = 36
$ws4.Range("F38").Value = 394
$ws4.Range("F41").Value = 14
$ws4.Range("F42").Value = 294
$ws4.Range("F43").Value = 728
$ws4.Range("F45").Value = 57
$ws4.Range("F46").Value = 70
